$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(88, 8).Value = 2009.2  # H88 was 2025.9
$ws.Cells.Item(88, 10).Value = 2520  # J88 was 2547.8333
$ws.Cells.Item(88, 12).Value = 2520  # L88 was 2547.8333
$ws.Cells.Item(88, 14).Value = -3332  # N88 was -3359.8333

$ws.Cells.Item(91, 8).Value = 2009.2  # H91 was 2025.9
$ws.Cells.Item(91, 10).Value = 2520  # J91 was 2547.8333
$ws.Cells.Item(91, 12).Value = 2520  # L91 was 2547.8333
$ws.Cells.Item(91, 14).Value = -5328  # N91 was -5355.8333

$ws.Cells.Item(107, 8).Value = 369.47058  # H107 was 380.11765
$ws.Cells.Item(107, 9).Value = 377  # I107 was 388.3125
$ws.Cells.Item(107, 11).Value = 377  # K107 was 388.3125
$ws.Cells.Item(107, 13).Value = 1543  # M107 was 1531.6875

$ws.Cells.Item(113, 8).Value = 5117.5557  # H113 was 4955.3
$ws.Cells.Item(113, 10).Value = 6332.6665  # J113 was 5623.25
$ws.Cells.Item(113, 12).Value = 6332.6665  # L113 was 5623.25
$ws.Cells.Item(113, 14).Value = -12840.6665  # N113 was -12131.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 16009.777  # H32 was 15909.437
$ws.Cells.Item(32, 9).Value = 6672.548  # I32 was 6716.976
$ws.Cells.Item(32, 10).Value = 29081.9  # J32 was 29222.654
$ws.Cells.Item(32, 11).Value = 6672.548  # K32 was 6716.976
$ws.Cells.Item(32, 12).Value = 29081.9  # L32 was 29222.654
$ws.Cells.Item(32, 13).Value = -6385.548  # M32 was -6429.976
$ws.Cells.Item(32, 14).Value = -29655.9  # N32 was -29796.654

$ws.Cells.Item(74, 8).Value = 3361.5667  # H74 was 3556.25
$ws.Cells.Item(74, 9).Value = 1338.0526  # I74 was 1420.6471
$ws.Cells.Item(74, 11).Value = 1338.0526  # K74 was 1420.6471
$ws.Cells.Item(74, 13).Value = -464.0526  # M74 was -546.6470999999999

$ws.Cells.Item(77, 8).Value = 3361.5667  # H77 was 3556.25
$ws.Cells.Item(77, 9).Value = 1338.0526  # I77 was 1420.6471
$ws.Cells.Item(77, 11).Value = 6690.263  # K77 was 7103.2355
$ws.Cells.Item(77, 13).Value = -2322.263  # M77 was -2735.2355

$ws.Cells.Item(122, 8).Value = 297165.47  # H122 was 288686.9
$ws.Cells.Item(122, 9).Value = 419110.38  # I122 was 387040.53
$ws.Cells.Item(122, 10).Value = 4497.7  # J122 was 4554.222
$ws.Cells.Item(122, 11).Value = 1257331.14  # K122 was 1161121.59
$ws.Cells.Item(122, 12).Value = 13493.1  # L122 was 13662.666
$ws.Cells.Item(122, 13).Value = -1254881.14  # M122 was -1158671.59
$ws.Cells.Item(122, 14).Value = -18393.1  # N122 was -18562.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 3411.238  # H134 was 3310.7273
$ws.Cells.Item(134, 9).Value = 2626.1177  # I134 was 2546.889
$ws.Cells.Item(134, 11).Value = 7878.353099999999  # K134 was 7640.667
$ws.Cells.Item(134, 13).Value = -5343.353099999999  # M134 was -5105.667

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 3939  # H16 was 4665
$ws.Cells.Item(16, 9).Value = 4087.5  # I16 was 4500
$ws.Cells.Item(16, 10).Value = 3345  # J16 was 4995
$ws.Cells.Item(16, 11).Value = 4087.5  # K16 was 4500
$ws.Cells.Item(16, 12).Value = 3345  # L16 was 4995
$ws.Cells.Item(16, 13).Value = -3800.5  # M16 was -4213
$ws.Cells.Item(16, 14).Value = -3919  # N16 was -5569

$ws.Cells.Item(107, 8).Value = 481.8  # H107 was 482.4
$ws.Cells.Item(107, 9).Value = 396.44446  # I107 was 433.5
$ws.Cells.Item(107, 10).Value = 1250  # J107 was 678
$ws.Cells.Item(107, 11).Value = 396.44446  # K107 was 433.5
$ws.Cells.Item(107, 12).Value = 1250  # L107 was 678
$ws.Cells.Item(107, 13).Value = 1523.55554  # M107 was 1486.5
$ws.Cells.Item(107, 14).Value = -5090  # N107 was -4518

$ws.Cells.Item(113, 8).Value = 3939  # H113 was 4665
$ws.Cells.Item(113, 9).Value = 4087.5  # I113 was 4500
$ws.Cells.Item(113, 10).Value = 3345  # J113 was 4995
$ws.Cells.Item(113, 11).Value = 4087.5  # K113 was 4500
$ws.Cells.Item(113, 12).Value = 3345  # L113 was 4995
$ws.Cells.Item(113, 13).Value = -1917.5  # M113 was -2330
$ws.Cells.Item(113, 14).Value = -7685  # N113 was -9335

$ws.Cells.Item(122, 8).Value = 910.5714  # H122 was 884.7
$ws.Cells.Item(122, 9).Value = 918.75  # I122 was 878.2857
$ws.Cells.Item(122, 11).Value = 2756.25  # K122 was 2634.8571
$ws.Cells.Item(122, 13).Value = -306.25  # M122 was -184.8571000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(52, 8).Value = 1150  # H52 was 1033.3334
$ws.Cells.Item(52, 10).Value = 1150  # J52 was 1033.3334
$ws.Cells.Item(52, 12).Value = 3450  # L52 was 3100.0002
$ws.Cells.Item(52, 14).Value = -3982  # N52 was -3632.0002

$ws.Cells.Item(55, 8).Value = 1981.1111  # H55 was 1228
$ws.Cells.Item(55, 10).Value = 2261.4285  # J55 was 1285
$ws.Cells.Item(55, 12).Value = 6784.2855  # L55 was 3855
$ws.Cells.Item(55, 14).Value = -7138.2855  # N55 was -4209

$ws.Cells.Item(60, 8).Value = 1148.5  # H60 was 854.8333
$ws.Cells.Item(60, 9).Value = 1598  # I60 was 858.25
$ws.Cells.Item(60, 10).Value = 699  # J60 was 848
$ws.Cells.Item(60, 11).Value = 4794  # K60 was 2574.75
$ws.Cells.Item(60, 12).Value = 2097  # L60 was 2544
$ws.Cells.Item(60, 13).Value = -4543  # M60 was -2323.75
$ws.Cells.Item(60, 14).Value = -2599  # N60 was -3046

$ws.Cells.Item(62, 8).Value = 8500  # H62 was 5285.7144
$ws.Cells.Item(62, 10).Value = 0  # J62 was 4000
$ws.Cells.Item(62, 12).Value = 0  # L62 was 12000
$ws.Cells.Item(62, 14).Value = ""  # N62 removed (was -13372)

$ws.Cells.Item(65, 8).Value = 8500  # H65 was 5285.7144
$ws.Cells.Item(65, 10).Value = 0  # J65 was 4000
$ws.Cells.Item(65, 12).Value = 0  # L65 was 36000
$ws.Cells.Item(65, 14).Value = ""  # N65 removed (was -42864)

$ws.Cells.Item(80, 8).Value = 5996.5  # H80 was 5996.857
$ws.Cells.Item(80, 9).Value = 5993  # I80 was 5994.5
$ws.Cells.Item(80, 11).Value = 17979  # K80 was 17983.5
$ws.Cells.Item(80, 13).Value = -17043  # M80 was -17047.5

$ws.Cells.Item(83, 8).Value = 5996.5  # H83 was 5996.857
$ws.Cells.Item(83, 9).Value = 5993  # I83 was 5994.5
$ws.Cells.Item(83, 11).Value = 53937  # K83 was 53950.5
$ws.Cells.Item(83, 13).Value = -49257  # M83 was -49270.5

$ws.Cells.Item(100, 8).Value = 450  # H100 was 0
$ws.Cells.Item(100, 10).Value = 450  # J100 was 0
$ws.Cells.Item(100, 12).Value = 1350  # L100 was 0
$ws.Cells.Item(100, 14).Value = -2972  # N100 was None

$ws.Cells.Item(102, 8).Value = 4600  # H102 was 3966.6667

$ws.Cells.Item(107, 8).Value = 2869.5715  # H107 was 2476.2222
$ws.Cells.Item(107, 9).Value = 5499  # I107 was 3832.6667
$ws.Cells.Item(107, 10).Value = 1817.8  # J107 was 1798
$ws.Cells.Item(107, 11).Value = 16497  # K107 was 11498.0001
$ws.Cells.Item(107, 12).Value = 5453.4  # L107 was 5394
$ws.Cells.Item(107, 13).Value = -14577  # M107 was -9578.000100000001
$ws.Cells.Item(107, 14).Value = -9293.4  # N107 was -9234

$ws.Cells.Item(108, 8).Value = 4361.7  # H108 was 4423.7
$ws.Cells.Item(108, 9).Value = 602.8333  # I108 was 706.1667
$ws.Cells.Item(108, 11).Value = 1808.4999  # K108 was 2118.5001
$ws.Cells.Item(108, 13).Value = 1071.5001  # M108 was 761.4998999999998

$ws.Cells.Item(109, 8).Value = 1400  # H109 was 999
$ws.Cells.Item(109, 9).Value = 1400  # I109 was 999
$ws.Cells.Item(109, 11).Value = 4200  # K109 was 2997
$ws.Cells.Item(109, 13).Value = -3160  # M109 was -1957

$ws.Cells.Item(111, 8).Value = 1513.5  # H111 was 1256.75
$ws.Cells.Item(111, 9).Value = 1513.5  # I111 was 1256.75
$ws.Cells.Item(111, 11).Value = 4540.5  # K111 was 3770.25
$ws.Cells.Item(111, 13).Value = -1473.5  # M111 was -703.25

$ws.Cells.Item(112, 8).Value = 3300  # H112 was 3750
$ws.Cells.Item(112, 9).Value = 0  # I112 was 2500
$ws.Cells.Item(112, 10).Value = 3300  # J112 was 5000
$ws.Cells.Item(112, 11).Value = 0  # K112 was 7500
$ws.Cells.Item(112, 13).Value = 9900  # M112 was -6392
$ws.Cells.Item(112, 14).Value = -12116  # N112 was -17216

$ws.Cells.Item(129, 8).Value = 2371.6  # H129 was 2183.182
$ws.Cells.Item(129, 10).Value = 2056  # J129 was 1805
$ws.Cells.Item(129, 12).Value = 6168  # L129 was 5415
$ws.Cells.Item(129, 14).Value = -16168  # N129 was -15415

$ws.Cells.Item(132, 8).Value = 2994  # H132 was 2047
$ws.Cells.Item(132, 10).Value = 2994  # J132 was 2047
$ws.Cells.Item(132, 12).Value = 26946  # L132 was 18423
$ws.Cells.Item(132, 14).Value = -32006  # N132 was -23483

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(46, 10).Value = 20000  # J46 was 0
$ws.Cells.Item(46, 12).Value = 20000  # L46 was 0
$ws.Cells.Item(46, 14).Value = -20312  # N46 newly added

$ws.Cells.Item(80, 9).Value = 598.3333  # I80 was 650
$ws.Cells.Item(80, 10).Value = 9382.5  # J80 was 6420
$ws.Cells.Item(80, 11).Value = 598.3333  # K80 was 650
$ws.Cells.Item(80, 12).Value = 9382.5  # L80 was 6420
$ws.Cells.Item(80, 13).Value = 399.6667  # M80 was 348
$ws.Cells.Item(80, 14).Value = -11378.5  # N80 was -8416

$ws.Cells.Item(83, 9).Value = 598.3333  # I83 was 650
$ws.Cells.Item(83, 10).Value = 9382.5  # J83 was 6420
$ws.Cells.Item(83, 11).Value = 2991.6665  # K83 was 3250
$ws.Cells.Item(83, 12).Value = 46912.5  # L83 was 32100
$ws.Cells.Item(83, 13).Value = 2000.3335  # M83 was 1742
$ws.Cells.Item(83, 14).Value = -56896.5  # N83 was -42084

$ws.Cells.Item(102, 8).Value = 2483.3  # H102 was 2561.842
$ws.Cells.Item(102, 9).Value = 1689.6923  # I102 was 1747.9166
$ws.Cells.Item(102, 11).Value = 1689.6923  # K102 was 1747.9166
$ws.Cells.Item(102, 13).Value = -67.69229999999993  # M102 was -125.9166

$ws.Cells.Item(122, 8).Value = 527868.4399999999  # H122 was 462176.5
$ws.Cells.Item(122, 9).Value = 94531.63  # I122 was 74774.78999999999
$ws.Cells.Item(122, 11).Value = 283594.89  # K122 was 224324.37
$ws.Cells.Item(122, 13).Value = -281144.89  # M122 was -221874.37

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 3779.5557  # H61 was 3564.9656
$ws.Cells.Item(61, 9).Value = 3480.5652  # I61 was 3255.56
$ws.Cells.Item(61, 11).Value = 3480.5652  # K61 was 3255.56
$ws.Cells.Item(61, 13).Value = -3278.5652  # M61 was -3053.56

$ws.Cells.Item(100, 8).Value = 1797.1111  # H100 was 2157.6
$ws.Cells.Item(100, 9).Value = 1771.75  # I100 was 2197
$ws.Cells.Item(100, 11).Value = 1771.75  # K100 was 2197
$ws.Cells.Item(100, 13).Value = -1230.75  # M100 was -1656

$ws.Cells.Item(113, 8).Value = 3779.5557  # H113 was 3564.9656
$ws.Cells.Item(113, 9).Value = 3480.5652  # I113 was 3255.56
$ws.Cells.Item(113, 11).Value = 3480.5652  # K113 was 3255.56
$ws.Cells.Item(113, 13).Value = -1310.5652  # M113 was -1085.56

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 1974.1666  # H122 was 2131.1667
$ws.Cells.Item(122, 9).Value = 1877.8889  # I122 was 1949.5
$ws.Cells.Item(122, 10).Value = 2263  # J122 was 2494.5
$ws.Cells.Item(122, 11).Value = 5633.6667  # K122 was 5848.5
$ws.Cells.Item(122, 12).Value = 6789  # L122 was 7483.5
$ws.Cells.Item(122, 13).Value = -3183.6667  # M122 was -3398.5
$ws.Cells.Item(122, 14).Value = -11689  # N122 was -12383.5
